# Linear regression workbook: add a "Tabelle2" worksheet with a mean-deviation
# (x-xbar)/(y-ybar) scratch calculation, wire up the `xbar`/`ybar` defined
# names it needs, and move the active-tab/selection onto the new sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- sheet1: selection moves from F10 to the B2:C7 data block -------------
$ws1.Range("B2:C7").Select()

# --- new sheet, inserted right after Tabelle1 ------------------------------
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "Tabelle2"

# match the metric (2cm) top/bottom page margins used on Tabelle1
$ws2.PageSetup.TopMargin = 56.692913399999995
$ws2.PageSetup.BottomMargin = 56.692913399999995

# --- workbook-level defined names used by the new sheet's formulas --------
$wb.Names.Add("xbar", '=Tabelle2!$B$11')
$wb.Names.Add("ybar", '=Tabelle2!$C$11')

# --- header row -------------------------------------------------------------
$ws2.Range("B3").Value = "x"
$ws2.Range("C3").Value = "y"
$ws2.Range("D3").Value = "x-xbar"
$ws2.Range("E3").Value = "y-ybar"
$ws2.Range("F3").Value = "(x-xbar)(y-ybar)"
$ws2.Range("G3").Value = "(x-xbar)^2"

# --- raw x/y sample plus the per-row deviation calculations ----------------
$xs = @(1,2,3,4,5,6)
$ys = @(1,3,2,5,4,5)
for ($i = 0; $i -lt 6; $i++) {
    $r = 4 + $i
    $ws2.Cells.Item($r,2).Value = $xs[$i]
    $ws2.Cells.Item($r,3).Value = $ys[$i]
    $ws2.Cells.Item($r,4).Formula = "=B$r-xbar"
    $ws2.Cells.Item($r,5).Formula = "=C$r-ybar"
    $ws2.Cells.Item($r,6).Formula = "=D$r*E$r"
    $ws2.Cells.Item($r,7).Formula = "=D$r^2"
}

# --- totals / averages -------------------------------------------------------
$ws2.Range("B10").Formula = "=SUM(B4:B9)"
$ws2.Range("C10").Formula = "=SUM(C4:C9)"
$ws2.Range("F10").Formula = "=SUM(F4:F9)"
$ws2.Range("G10").Formula = "=SUM(G4:G9)"

$ws2.Range("B11").Formula = "=AVERAGE(B4:B9)"
$ws2.Range("C11").Formula = "=AVERAGE(C4:C9)"

# --- slope / intercept / check ----------------------------------------------
$ws2.Range("F13").Formula = "=F10/G10"
$ws2.Range("C15").Formula = "=7*F13+F16"
$ws2.Range("F16").Formula = "=ybar-F13*xbar"

# --- view: Tabelle2 becomes the active/selected sheet -----------------------
$ws2.Range("B3").Select()
$ws2.Activate()

$wb.Application.Calculate()
